$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/26_kapako1.wav"
$ws.Range("B2").Value = "pngimages/26_pineapple.png"

$ws.Range("A3").Value = "trainingaudio/20_tatito1.wav"
$ws.Range("B3").Value = "pngimages/20_pizza.png"

$ws.Range("A4").Value = "trainingaudio/06_titoka3.wav"
$ws.Range("B4").Value = "pngimages/06_tent.png"

$ws.Range("A5").Value = "trainingaudio/04_kitoti2.wav"
$ws.Range("B5").Value = "pngimages/04_ladder.png"

$ws.Range("A6").Value = "trainingaudio/08_tipako2.wav"
$ws.Range("B6").Value = "pngimages/08_bell.png"

$ws.Range("A7").Value = "trainingaudio/12_pokika3.wav"
$ws.Range("B7").Value = "pngimages/12_pie.png"
